$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new column at DG (shifts DG:EK -> DH:EL) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("DG1").EntireColumn.Insert()

# Header for the new column
$wsPrix.Range("DG1").Value = "02-nov"

# Fill new column with "-" placeholder for data rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 111).Value = "-"
}

# --- Sheet "Gaz": append new row 139 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$cA = $wsGaz.Cells.Item(139, 1)
$cA.NumberFormat = "@"
$cA.Value = "2025-10-31"
$cA.ClearFormats()
$wsGaz.Cells.Item(139, 2).Value = 29.9

# --- Sheet "CO2": append new row 139 ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$cA2 = $wsCo2.Cells.Item(139, 1)
$cA2.NumberFormat = "@"
$cA2.Value = "2025-10-31"
$cA2.ClearFormats()
$wsCo2.Cells.Item(139, 2).Value = 78
